$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "implemented"
$ws.Range("E6").Value = "Risk parity backtest implemented (ERC) with rolling window + bounds + UI presets + tests/help."
$ws.Range("F6").Value = "27/12/2025 03:49"
